$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 159, shifting existing rows 159-167 down to 160-168
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new data record
$ws.Cells.Item(159, 1).Value2 = 10
$ws.Cells.Item(159, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(159, 3).Value2 = "La Araucanía"
$ws.Cells.Item(159, 4).Value2 = 44509
$ws.Cells.Item(159, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(159, 5).Value2 = 9
$ws.Cells.Item(159, 6).Value2 = "Fruta"
$ws.Cells.Item(159, 7).Value2 = 100101
$ws.Cells.Item(159, 8).Value2 = "Berries"
$ws.Cells.Item(159, 9).Value2 = 100112025
$ws.Cells.Item(159, 10).Value2 = "Frutilla"
$ws.Cells.Item(159, 11).Value2 = "Sin especificar"
$ws.Cells.Item(159, 12).Value2 = "Primera"
$ws.Cells.Item(159, 13).Value2 = 270
$ws.Cells.Item(159, 14).Value2 = 9000
$ws.Cells.Item(159, 15).Value2 = 9000
$ws.Cells.Item(159, 16).Value2 = 9000
$ws.Cells.Item(159, 17).Value2 = "$/bandeja 7 kilos"
$ws.Cells.Item(159, 18).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(159, 19).Value2 = 1286
$ws.Cells.Item(159, 20).Value2 = 7
